# river update May 2024
# Adds three new rows (30-32) of river-monitoring results for
# "Manganui o te Ao at Ashworth" covering the 2019 - 2023 period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: ASPM
$ws.Cells.Item(30, 1).Value = "Manganui o te Ao at Ashworth"
$ws.Cells.Item(30, 2).Value = "ASPM"
$ws.Cells.Item(30, 3).Value = "C"
$ws.Cells.Item(30, 4).Value = "2019 - 2023"
$ws.Cells.Item(30, 5).Value = "RepSite"
$ws.Cells.Item(30, 6).Value = 0.358
$ws.Cells.Item(30, 7).Value = 0.3342
$ws.Cells.Item(30, 8).Value = 0.391
$ws.Cells.Item(30, 9).Value = 0.391
$ws.Cells.Item(30, 12).Value = 0.314
$ws.Cells.Item(30, 13).Value = 0.39065
$ws.Cells.Item(30, 14).Value = 0.391
$ws.Cells.Item(30, 15).Value = 1789685
$ws.Cells.Item(30, 16).Value = 5646155
$ws.Cells.Item(30, 17).Value = "Ruapehu District"
$ws.Cells.Item(30, 18).Value = "Whanganui"
$ws.Cells.Item(30, 19).Value = "Pipiriki"
$ws.Cells.Item(30, 20).Value = "Whai_5i"

# Row 31: MCI
$ws.Cells.Item(31, 1).Value = "Manganui o te Ao at Ashworth"
$ws.Cells.Item(31, 2).Value = "MCI"
$ws.Cells.Item(31, 3).Value = "C"
$ws.Cells.Item(31, 4).Value = "2019 - 2023"
$ws.Cells.Item(31, 5).Value = "RepSite"
$ws.Cells.Item(31, 6).Value = 104.76
$ws.Cells.Item(31, 7).Value = 101.672
$ws.Cells.Item(31, 8).Value = 112
$ws.Cells.Item(31, 9).Value = 112
$ws.Cells.Item(31, 12).Value = 104.3
$ws.Cells.Item(31, 13).Value = 111.16
$ws.Cells.Item(31, 14).Value = 112
$ws.Cells.Item(31, 15).Value = 1789685
$ws.Cells.Item(31, 16).Value = 5646155
$ws.Cells.Item(31, 17).Value = "Ruapehu District"
$ws.Cells.Item(31, 18).Value = "Whanganui"
$ws.Cells.Item(31, 19).Value = "Pipiriki"
$ws.Cells.Item(31, 20).Value = "Whai_5i"

# Row 32: QMCI
$ws.Cells.Item(32, 1).Value = "Manganui o te Ao at Ashworth"
$ws.Cells.Item(32, 2).Value = "QMCI"
$ws.Cells.Item(32, 3).Value = "D"
$ws.Cells.Item(32, 4).Value = "2019 - 2023"
$ws.Cells.Item(32, 5).Value = "RepSite"
$ws.Cells.Item(32, 6).Value = 3.69
$ws.Cells.Item(32, 7).Value = 3.6614
$ws.Cells.Item(32, 8).Value = 4.29
$ws.Cells.Item(32, 9).Value = 4.29
$ws.Cells.Item(32, 12).Value = 3.494
$ws.Cells.Item(32, 13).Value = 4.21615
$ws.Cells.Item(32, 14).Value = 4.29
$ws.Cells.Item(32, 15).Value = 1789685
$ws.Cells.Item(32, 16).Value = 5646155
$ws.Cells.Item(32, 17).Value = "Ruapehu District"
$ws.Cells.Item(32, 18).Value = "Whanganui"
$ws.Cells.Item(32, 19).Value = "Pipiriki"
$ws.Cells.Item(32, 20).Value = "Whai_5i"
